$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, copying the exact formatting of the other header cells
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Add time_taken values for each data row (as plain text)
$timestamps = @(
    "2021-10-05 13:41:01.561445",
    "2021-10-05 13:41:01.561457",
    "2021-10-05 13:41:01.561461",
    "2021-10-05 13:41:01.561464",
    "2021-10-05 13:41:01.561467",
    "2021-10-05 13:41:01.561470",
    "2021-10-05 13:41:01.561474",
    "2021-10-05 13:41:01.561477",
    "2021-10-05 13:41:01.561481",
    "2021-10-05 13:41:01.561484",
    "2021-10-05 13:41:01.561487",
    "2021-10-05 13:41:01.561490",
    "2021-10-05 13:41:01.561493",
    "2021-10-05 13:41:01.561496",
    "2021-10-05 13:41:01.561499",
    "2021-10-05 13:41:01.561502",
    "2021-10-05 13:41:01.561506",
    "2021-10-05 13:41:01.561509",
    "2021-10-05 13:41:01.561512"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

$wb.Save()
